# Adds a "2021" data column (Policia, OPM and DCR tabs) before the
# existing "total" column, shifting the total column from F to G and
# recomputing the totals to include the new 2021 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "total" column (F) one column to the right (G),
# preserving its header text.
$ws.Range("G1").Value = $ws.Range("F1").Value2

# New "2021" values that populate the inserted column F.
$year2021 = @{
    2 = 29
    3 = 28
    4 = 7
    5 = 2
    6 = 0
    7 = 3
    8 = 0
    9 = 69
    10 = 69
}

# Recomputed totals (now covering 2017-2021) that move into column G.
$totals = @{
    2 = 192
    3 = 138
    4 = 45
    5 = 26
    6 = 1
    7 = 23
    8 = 3
    9 = 428
    10 = 434
}

$ws.Range("F1").Value = 2021

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = $year2021[$r]
    $ws.Cells.Item($r, 7).Value = $totals[$r]
}

# Apply the new "Aptos Narrow" 11pt font to the new 2021 column (F1:F10)
# and to the data rows of the shifted total column (G2:G10), matching
# the formatting applied by Excel's font picker. Building it once on a
# scratch cell and pasting the format across avoids creating a separate
# style for every destination cell.
$scratch = $ws.Range("Z1")
$scratch.Font.Size = 11
$scratch.Font.Name = "Aptos Narrow"
$scratch.Copy()

$ws.Range("F1:F10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G2:G10").PasteSpecial(-4122)  # xlPasteFormats

$scratch.Clear()
$excel.CutCopyMode = $false
